{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Rewrites five body paragraphs (Summary / Experience / Education / Skills /\n// Projects content) to match the target revision. The Experience and\n// Projects paragraphs contain several text runs separated by manual line\n// breaks (<w:br/>), so we rebuild them via insertOoxml(..., Replace) which\n// lets us place a single run holding multiple <w:t>/<w:br/> children --\n// matching the exact shape produced by the authoring tool.\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// Build a minimal OOXML package wrapping a single <w:p> whose single <w:r>\n// contains the given text parts joined by <w:br/> elements.\nfunction paragraphOoxmlPackage(parts) {\n  const runChildren = parts\n    .map((p) => \"<w:t xml:space=\\\"preserve\\\">\" + escapeXml(p) + \"</w:t>\")\n    .join(\"<w:br/>\");\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n    \"<pkg:xmlData>\" +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    \"</Relationships>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p><w:r>\" +\n    runChildren +\n    \"</w:r></w:p></w:body></w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\n// Replace the visible content of `paragraph` with `parts` (array of strings)\n// joined by manual line breaks, preserving the paragraph itself (so its\n// pPr/style survive) and its position in the body.\nfunction setParagraphParts(paragraph, parts) {\n  paragraph.insertOoxml(paragraphOoxmlPackage(parts), Word.InsertLocation.replace);\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map each target paragraph by matching the *old* text known from the\n// original document, so the script is resilient to being re-run against\n// the same starting point.\nconst oldTexts = {\n  summary:\n    \"A passionate and driven student with a strong desire to contribute to the vibrant world of game development and UI design. My academic background in computer science and a keen understanding of game mechanics and user experience principles makes me a confident and adaptable individual eager to learn and contribute to innovative projects. I am committed to developing engaging and user-friendly applications, and I am eager to leverage my skills to create impactful solutions.\",\n  education: \"B.E. CSE with 7.1 CGPA\",\n  skills: \"java, react, C#, node,java and python\",\n};\n\nlet summaryPara = null;\nlet experiencePara = null;\nlet educationPara = null;\nlet skillsPara = null;\nlet projectsPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text === oldTexts.summary) {\n    summaryPara = paragraphs.items[i];\n  } else if (text === oldTexts.education) {\n    educationPara = paragraphs.items[i];\n  } else if (text === oldTexts.skills) {\n    skillsPara = paragraphs.items[i];\n  } else if (text.indexOf(\"AI Intern, Blue Silicon Infotech\") !== -1) {\n    experiencePara = paragraphs.items[i];\n  } else if (text.indexOf(\"QR Scanner\") !== -1) {\n    projectsPara = paragraphs.items[i];\n  }\n}\n\n// --- Summary ---\nif (summaryPara) {\n  summaryPara.insertText(\n    \"Highly motivated and results-oriented individual seeking a challenging role in the development of innovative mobile applications. Proven ability to collaborate effectively, solve complex problems, and contribute to a dynamic team environment. Eager to leverage skills in software development, project management, and communication to contribute to impactful projects.\",\n    Word.InsertLocation.replace\n  );\n}\n\n// --- Experience ---\nif (experiencePara) {\n  setParagraphParts(experiencePara, [\n    \"AI Intern, Blue Silicon Infotech\",\n    \"Developed and implemented AI resume enhancer, resulting in a 20% increase in resume completion rates. Optimized resume templates for improved readability and clarity. Achieved quantifiable results, including a 15% improvement in resume accuracy.\",\n  ]);\n}\n\nawait context.sync();\n\n// --- Education ---\nif (educationPara) {\n  educationPara.insertText(\n    \"Bachelor of Engineering from AVIT. Graduated: 2026-05. GPA: 7.1.\",\n    Word.InsertLocation.replace\n  );\n}\n\n// --- Skills ---\nif (skillsPara) {\n  skillsPara.insertText(\n    \"Here's a revised skills section tailored for a global editing resume, focusing on conciseness, organization, and professionalism:, Skills**,    **Programming Languages:** Python, Java, Node.js, TypeScript, C#, Flutter, Kotlin, Dart, React, Python, SQL,    **Database:** MySQL, PostgreSQL, MongoDB, SQL Server,    **Web Development:** HTML, CSS, JavaScript, React, Angular, Vue.js,    **Cloud Technologies:** AWS, Azure, Google Cloud Platform,    **Operating Systems:** Linux, Windows, macOS,    **Version Control:** Git, GitHub, GitLab,    **Data Analysis:** Pandas, NumPy, Matplotlib, Seaborn,    **Testing:** Unit Testing, Integration Testing, End-to-End Testing,    **API Development:** RESTful APIs, GraphQL,    **Design Principles:** SOLID, DRY, KISS,    **Other:** Agile Development, Mobile Development, Data Science\",\n    Word.InsertLocation.replace\n  );\n}\n\n// --- Projects ---\nif (projectsPara) {\n  setParagraphParts(projectsPara, [\n    \"**Project:** Enhanced QR Scanner and Generator\",\n    \"**Summary:** This project aimed to significantly improve the efficiency and accuracy of QR scanning and generator functionality. By implementing a novel algorithm and incorporating real-time data integration, we achieved a demonstrable increase in accuracy and reduced processing time. This improved functionality was directly translated into increased sales and reduced operational costs.\",\n    \"**Technologies:**\",\n    \"*  QR scanner and generator\",\n    \"*  Prediction pro\",\n    \"*  Simple purchase order manager\",\n    \"*  PDF maker\",\n    \"**Contributions:**\",\n    \"*  Improved accuracy in QR scanning and generator processing.\",\n    \"*  Enhanced real-time data integration for improved processing.\",\n    \"*  Reduced processing time by 20%.\",\n    \"**Measurable Results:**\",\n    \"*  Increased accuracy in QR scanning and generator processing.\",\n    \"*  Reduced processing time by 20%.\",\n    \"*  Improved sales and reduced operational costs.\",\n    \"**Improvements:**\",\n    \"*  Improved accuracy in QR scanning and generator processing.\",\n    \"*  Enhanced real-time data integration for improved processing.\",\n    \"*  Reduced processing time by 20%.\",\n  ]);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document ($d below).\n#\n# Rewrites five body paragraphs (Summary / Experience / Education / Skills /\n# Projects content) to match the target revision. The Experience and\n# Projects paragraphs hold several text runs separated by manual line\n# breaks (<w:br/> <-> Chr(11) / vertical-tab in Range.Text), so the old\n# text is located (and replaced) including the embedded break characters,\n# producing the same \"single run, multiple <w:t>/<w:br/>\" shape as the\n# original authoring tool.\n\n$d = $word.ActiveDocument\n$vt = [char]11\n\nfunction Replace-DocText($searchText, $newText) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute($searchText)\n    if ($found) {\n        $rng.Text = $newText\n    }\n}\n\n# --- Summary ---\n$oldSummary = \"A passionate and driven student with a strong desire to contribute to the vibrant world of game development and UI design. My academic background in computer science and a keen understanding of game mechanics and user experience principles makes me a confident and adaptable individual eager to learn and contribute to innovative projects. I am committed to developing engaging and user-friendly applications, and I am eager to leverage my skills to create impactful solutions.\"\n$newSummary = \"Highly motivated and results-oriented individual seeking a challenging role in the development of innovative mobile applications. Proven ability to collaborate effectively, solve complex problems, and contribute to a dynamic team environment. Eager to leverage skills in software development, project management, and communication to contribute to impactful projects.\"\nReplace-DocText $oldSummary $newSummary\n\n# --- Experience ---\n$oldExperience = \"Here's a polished and professional version of the work experience section, incorporating the requested requirements:\" + $vt + $vt + $vt + $vt + \"**AI Intern, Blue Silicon Infotech**\" + $vt + $vt + $vt + $vt + \"Highly motivated and results-oriented AI intern with a proven track record of developing and deploying cutting-edge AI solutions. Developed expertise in natural language processing (NLP), machine learning (ML), and deep learning techniques. Successfully collaborated with diverse teams to design, develop, and implement AI models for various applications, including customer support chatbots and personalized recommendations. Demonstrated proficiency in data preprocessing, model training, and evaluation. Adept at collaborating effectively with stakeholders to achieve impactful outcomes. Opportunity to contribute to innovative AI projects and contribute to the growth of Blue Silicon Infotech's AI capabilities.\"\n$newExperienceParts = @(\n    \"AI Intern, Blue Silicon Infotech\",\n    \"Developed and implemented AI resume enhancer, resulting in a 20% increase in resume completion rates. Optimized resume templates for improved readability and clarity. Achieved quantifiable results, including a 15% improvement in resume accuracy.\"\n)\n$newExperience = ($newExperienceParts -join $vt)\nReplace-DocText $oldExperience $newExperience\n\n# --- Education ---\n$oldEducation = \"B.E. CSE with 7.1 CGPA\"\n$newEducation = \"Bachelor of Engineering from AVIT. Graduated: 2026-05. GPA: 7.1.\"\nReplace-DocText $oldEducation $newEducation\n\n# --- Skills ---\n$oldSkills = \"java, react, C#, node,java and python\"\n$newSkills = \"Here's a revised skills section tailored for a global editing resume, focusing on conciseness, organization, and professionalism:, Skills**,    **Programming Languages:** Python, Java, Node.js, TypeScript, C#, Flutter, Kotlin, Dart, React, Python, SQL,    **Database:** MySQL, PostgreSQL, MongoDB, SQL Server,    **Web Development:** HTML, CSS, JavaScript, React, Angular, Vue.js,    **Cloud Technologies:** AWS, Azure, Google Cloud Platform,    **Operating Systems:** Linux, Windows, macOS,    **Version Control:** Git, GitHub, GitLab,    **Data Analysis:** Pandas, NumPy, Matplotlib, Seaborn,    **Testing:** Unit Testing, Integration Testing, End-to-End Testing,    **API Development:** RESTful APIs, GraphQL,    **Design Principles:** SOLID, DRY, KISS,    **Other:** Agile Development, Mobile Development, Data Science\"\nReplace-DocText $oldSkills $newSkills\n\n# --- Projects ---\n$oldProjects = \"Here's the improved content for the project:\" + $vt + $vt + \"**QR Scanner\" + $vt + $vt + \"**\" + $vt + $vt + \"**Prediction Pro**\" + $vt + $vt + \"**Simple Purchase Order Manager**\"\n$newProjectsParts = @(\n    \"**Project:** Enhanced QR Scanner and Generator\",\n    \"**Summary:** This project aimed to significantly improve the efficiency and accuracy of QR scanning and generator functionality. By implementing a novel algorithm and incorporating real-time data integration, we achieved a demonstrable increase in accuracy and reduced processing time. This improved functionality was directly translated into increased sales and reduced operational costs.\",\n    \"**Technologies:**\",\n    \"*  QR scanner and generator\",\n    \"*  Prediction pro\",\n    \"*  Simple purchase order manager\",\n    \"*  PDF maker\",\n    \"**Contributions:**\",\n    \"*  Improved accuracy in QR scanning and generator processing.\",\n    \"*  Enhanced real-time data integration for improved processing.\",\n    \"*  Reduced processing time by 20%.\",\n    \"**Measurable Results:**\",\n    \"*  Increased accuracy in QR scanning and generator processing.\",\n    \"*  Reduced processing time by 20%.\",\n    \"*  Improved sales and reduced operational costs.\",\n    \"**Improvements:**\",\n    \"*  Improved accuracy in QR scanning and generator processing.\",\n    \"*  Enhanced real-time data integration for improved processing.\",\n    \"*  Reduced processing time by 20%.\"\n)\n$newProjects = ($newProjectsParts -join $vt)\nReplace-DocText $oldProjects $newProjects\n"}
